# Update the "Förändrad" date (column C) from 2023-09-17 (45186) to
# 2023-09-19 (45188) for every data row (rows 2-51) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
